$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.535.13"
$ws.Range("E2").Value = '  +0.12%  '
$ws.Range("D3").Value = "'1.737.16"
$ws.Range("E3").Value = '  +0.24%  '
$ws.Range("D4").Value = "'0.9998"
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").Value = "'246.93"
$ws.Range("E5").Value = '  +1.11%  '
$ws.Range("D6").Value = "'1.000"
$ws.Range("E6").Value = '  +0.02%  '
$ws.Range("D7").Value = "'0.4911"
$ws.Range("E7").Value = '  +2.46%  '
$ws.Range("D8").Value = "'0.2667"
$ws.Range("E8").Value = '  -0.17%  '
$ws.Range("D9").Value = "'0.06302"
$ws.Range("E9").Value = '  +1.20%  '
$ws.Range("D10").Value = "'1.733.02"
$ws.Range("E10").Value = '  -0.04%  '
$ws.Range("D11").Value = "'0.07054"
$ws.Range("E11").Value = '  -1.04%  '
$ws.Range("D12").Value = "'15.71"
$ws.Range("E12").Value = '  -0.13%  '
$ws.Range("D13").Value = "'4.603"
$ws.Range("E13").Value = '  +1.55%  '
$ws.Range("D14").Value = "'0.6115"
$ws.Range("E14").Value = '  -0.66%  '
$ws.Range("D15").Value = "'77.45"
$ws.Range("E15").Value = '  +0.71%  '
$ws.Range("D16").Value = "'0.9999"
$ws.Range("E16").Value = '  -0.05%  '
$ws.Range("D17").Value = "'0.000007398"
$ws.Range("E17").Value = '  +7.01%  '
$ws.Range("D18").Value = "'26.534.27"
$ws.Range("E18").Value = '  +0.17%  '
$ws.Range("D19").Value = "'1.000"
$ws.Range("E19").Value = '  -0.01%  '
$ws.Range("D20").Value = "'11.54"
$ws.Range("E20").Value = '  -1.58%  '
$ws.Range("D21").Value = "'1.957.22"
$ws.Range("E21").Value = '  -0.09%  '
$ws.Range("D22").Value = "'4.590"
$ws.Range("E22").Value = '  +0.33%  '
$ws.Range("D23").Value = "'8.719"
$ws.Range("E23").Value = '  -1.97%  '
$ws.Range("D24").Value = "'5.254"
$ws.Range("E24").Value = '  -1.43%  '
$ws.Range("D25").Value = "'140.14"
$ws.Range("E25").Value = '  +2.87%  '
$ws.Range("D26").Value = "'15.47"
$ws.Range("E26").Value = '  +0.65%  '
$ws.Range("D27").Value = "'1.419"
$ws.Range("E27").Value = '  +0.38%  '
$ws.Range("D28").Value = "'1.765"
$ws.Range("E28").Value = '  -1.71%  '
$ws.Range("D29").Value = "'107.88"
$ws.Range("E29").Value = '  +0.92%  '
$ws.Range("D30").Value = "'4.050"
$ws.Range("E30").Value = '  +1.61%  '
$ws.Range("D31").Value = "'0.08053"
$ws.Range("E31").Value = '  +0.80%  '
$ws.Range("D32").Value = "'3.720"
$ws.Range("E32").Value = '  +0.17%  '
$ws.Range("D33").Value = "'0.04592"
$ws.Range("E33").Value = '  +1.10%  '
$ws.Range("D34").Value = "'0.9998"
$ws.Range("E34").Value = '  -0.02%  '
$ws.Range("D35").Value = "'2.611"
$ws.Range("E35").Value = '  -0.22%  '
$ws.Range("D36").Value = "'1.009"
$ws.Range("E36").Value = '  +1.79%  '
$ws.Range("D37").Value = "'0.6372"
$ws.Range("E37").Value = '  +0.00%  '
$ws.Range("D38").Value = "'0.8965"
$ws.Range("E38").Value = '  -4.06%  '
$ws.Range("D39").Value = "'2.019"
$ws.Range("E39").Value = '  +1.87%  '
$ws.Range("D40").Value = "'2.404"
$ws.Range("E40").Value = '  -0.25%  '
$ws.Range("D41").Value = "'1.005"
$ws.Range("E41").Value = '  -0.08%  '
$ws.Range("D42").Value = "'0.01509"
$ws.Range("E42").Value = '  +0.10%  '
$ws.Range("D43").Value = "'102.37"
$ws.Range("E43").Value = '  -7.06%  '
$ws.Range("D44").Value = "'5.407"
$ws.Range("E44").Value = '  -5.12%  '
$ws.Range("D45").Value = "'0.3907"
$ws.Range("E45").Value = '  +0.13%  '
$ws.Range("D46").Value = "'6.893"
$ws.Range("E46").Value = '  -0.52%  '
$ws.Range("D47").Value = "'0.1187"
$ws.Range("E47").Value = '  -0.41%  '
$ws.Range("D48").Value = "'0.05399"
$ws.Range("E48").Value = '  +1.29%  '
$ws.Range("D49").Value = "'30.60"
$ws.Range("E49").Value = '  -0.53%  '
$ws.Range("D50").Value = "'7.821"
$ws.Range("E50").Value = '  -1.04%  '
$ws.Range("D51").Value = "'1.273"
$ws.Range("E51").Value = '  +0.34%  '
